# Nexial "#system" macro-catalog sheet update:
#  - add `outputToCloud(resource)` to the "base" command list
#  - add a brand new "text" command category (single entry: spellCheck(var,profile,text))
#
# Both the "target" (category names) list and the "base" (function names) list are
# kept in alphabetical order, so the new entries are inserted in-place (not appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$xlShiftToRight  = -4161

# 1) Insert a brand new column before the existing "web" column (column Y) to hold the
#    new "text" category. Everything from Y..AD shifts right to Z..AE.
$ws.Columns("Y:Y").Insert($xlShiftToRight)
$ws.Range("Y1").Value2 = "text"
$ws.Range("Y2").Value2 = "spellCheck(var,profile,text)"

# 2) Insert "text" into the "target" (category) list in column A, in alphabetical
#    order, right before "web" (currently row 25). Only column A's values need to
#    move, so shift them manually (bottom-up) instead of using a row/range Insert
#    (which would also disturb unrelated columns in this environment).
for ($r = 30; $r -ge 25; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(25, 1).Value2 = "text"

# 3) Insert "outputToCloud(resource)" into the "base" (function) list in column E, in
#    alphabetical order, right before "prependText(var,prependWith)" (currently row
#    22). Only column E's values need to move - shift manually as above.
for ($r = 38; $r -ge 22; $r--) {
    $ws.Cells.Item($r + 1, 5).Value2 = $ws.Cells.Item($r, 5).Value2
}
$ws.Cells.Item(22, 5).Value2 = "outputToCloud(resource)"

# 4) Fix up the named ranges so they cover the correct (now shifted/expanded) ranges.
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("base").RefersTo      = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
